$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$c = $ws.Range("A8").Characters(21, 2)
$c.Text = "25"
$c.Font.Size = 10
$c.Font.Name = "Andale WT"

$d1 = $ws.Range("C9").Characters(27, 9)
$d1.Text = "6/17/2024"
$d1.Font.Size = 10
$d1.Font.Name = "Andale WT"

$d2 = $ws.Range("C9").Characters(47, 9)
$d2.Text = "6/23/2024"
$d2.Font.Size = 10
$d2.Font.Name = "Andale WT"

# --- Weekly crime-statistics table updates (rows 16-31) ---
# C16: ('14', 's', '20') -> ('15', None, '1')
$ws.Range("I15").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$ws.Range("C16").Value = 1

# D16: ('14', 's', '20') -> ('15', None, '1')
$ws.Range("I15").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").Value = 1

# E16: ('14', 's', '21') -> ('16', None, '0')
$ws.Range("N20").Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null
$ws.Range("E16").Value = 0

# F16: ('15', None, '1') -> ('15', None, '2')
$ws.Range("F16").Value = 2

# H16: ('16', None, '-50') -> ('16', None, '0')
$ws.Range("H16").Value = 0

# I16: ('15', None, '21') -> ('15', None, '22')
$ws.Range("I16").Value = 22

# J16: ('15', None, '6') -> ('15', None, '7')
$ws.Range("J16").Value = 7

# K16: ('16', None, '250') -> ('16', None, '214.285714285714')
$ws.Range("K16").Value = 214.285714285714

# L16: ('16', None, '200') -> ('16', None, '214.285714285714')
$ws.Range("L16").Value = 214.285714285714

# M16: ('16', None, '133.333333333333') -> ('16', None, '57.142857142857')
$ws.Range("M16").Value = 57.142857142857

# N16: ('16', None, '-67.1875') -> ('16', None, '-69.863013698630')
$ws.Range("N16").Value = -69.863013698630

# C17: ('15', None, '2') -> ('14', 's', '20')
$ws.Range("C20").Copy() | Out-Null
$ws.Range("C17").PasteSpecial(-4122) | Out-Null
$ws.Range("C17").Value = "0"

# D17: ('15', None, '1') -> ('14', 's', '20')
$ws.Range("C20").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4122) | Out-Null
$ws.Range("D17").Value = "0"

# E17: ('16', None, '100') -> ('14', 's', '21')
$ws.Range("C20").Copy() | Out-Null
$ws.Range("E17").PasteSpecial(-4122) | Out-Null
$ws.Range("E17").Value = "***.*"

# F17: ('15', None, '3') -> ('15', None, '2')
$ws.Range("F17").Value = 2

# H17: ('16', None, '200') -> ('16', None, '100')
$ws.Range("H17").Value = 100

# N17: ('16', None, '-61.904761904761') -> ('16', None, '-63.636363636363')
$ws.Range("N17").Value = -63.636363636363

# G18: ('15', None, '2') -> ('15', None, '1')
$ws.Range("G18").Value = 1

# C19: ('14', 's', '20') -> ('15', None, '2')
$ws.Range("I15").Copy() | Out-Null
$ws.Range("C19").PasteSpecial(-4122) | Out-Null
$ws.Range("C19").Value = 2

# D19: ('15', None, '1') -> ('14', 's', '20')
$ws.Range("C20").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4122) | Out-Null
$ws.Range("D19").Value = "0"

# E19: ('16', None, '-100') -> ('14', 's', '21')
$ws.Range("C20").Copy() | Out-Null
$ws.Range("E19").PasteSpecial(-4122) | Out-Null
$ws.Range("E19").Value = "***.*"

# F19: ('15', None, '11') -> ('15', None, '5')
$ws.Range("F19").Value = 5

# G19: ('15', None, '6') -> ('15', None, '3')
$ws.Range("G19").Value = 3

# H19: ('16', None, '83.333333333333') -> ('16', None, '66.666666666666')
$ws.Range("H19").Value = 66.666666666666

# I19: ('15', None, '24') -> ('15', None, '26')
$ws.Range("I19").Value = 26

# K19: ('16', None, '50') -> ('16', None, '62.5')
$ws.Range("K19").Value = 62.5

# L19: ('16', None, '200') -> ('16', None, '188.888888888889')
$ws.Range("L19").Value = 188.888888888889

# M19: ('16', None, '-14.285714285714') -> ('16', None, '-16.129032258064')
$ws.Range("M19").Value = -16.129032258064

# N19: ('16', None, '-64.179104477611') -> ('16', None, '-64.864864864864')
$ws.Range("N19").Value = -64.864864864864

# C21: ('18', None, '2') -> ('18', None, '3')
$ws.Range("C21").Value = 3

# D21: ('18', None, '2') -> ('18', None, '1')
$ws.Range("D21").Value = 1

# E21: ('19', None, '0') -> ('19', None, '200')
$ws.Range("E21").Value = 200

# F21: ('18', None, '15') -> ('18', None, '9')
$ws.Range("F21").Value = 9

# G21: ('18', None, '11') -> ('18', None, '7')
$ws.Range("G21").Value = 7

# H21: ('19', None, '36.363636363636') -> ('19', None, '28.571428571428')
$ws.Range("H21").Value = 28.571428571428

# I21: ('18', None, '54') -> ('18', None, '57')
$ws.Range("I21").Value = 57

# J21: ('18', None, '29') -> ('18', None, '30')
$ws.Range("J21").Value = 30

# K21: ('19', None, '86.206896551724') -> ('19', None, '90')
$ws.Range("K21").Value = 90

# L21: ('19', None, '125') -> ('19', None, '128')
$ws.Range("L21").Value = 128

# M21: ('19', None, '28.571428571428') -> ('19', None, '14')
$ws.Range("M21").Value = 14

# N21: ('19', None, '-68.965517241379') -> ('19', None, '-70.157068062827')
$ws.Range("N21").Value = -70.157068062827

# C24: ('15', None, '1') -> ('14', 's', '20')
$ws.Range("C20").Copy() | Out-Null
$ws.Range("C24").PasteSpecial(-4122) | Out-Null
$ws.Range("C24").Value = "0"

# F24: ('15', None, '8') -> ('15', None, '5')
$ws.Range("F24").Value = 5

# G24: ('15', None, '2') -> ('15', None, '1')
$ws.Range("G24").Value = 1

# H24: ('16', None, '300') -> ('16', None, '400')
$ws.Range("H24").Value = 400

# L24: ('16', None, '25') -> ('16', None, '0')
$ws.Range("L24").Value = 0

# M24: ('16', None, '-46.428571428571') -> ('16', None, '-48.275862068965')
$ws.Range("M24").Value = -48.275862068965

# C26: ('15', None, '1') -> ('15', None, '2')
$ws.Range("C26").Value = 2

# D26: ('14', 's', '20') -> ('15', None, '3')
$ws.Range("I15").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null
$ws.Range("D26").Value = 3

# E26: ('14', 's', '21') -> ('16', None, '-33.333333333333')
$ws.Range("N20").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4122) | Out-Null
$ws.Range("E26").Value = -33.333333333333

# F26: ('15', None, '3') -> ('15', None, '4')
$ws.Range("F26").Value = 4

# G26: ('15', None, '8') -> ('15', None, '7')
$ws.Range("G26").Value = 7

# H26: ('16', None, '-62.5') -> ('16', None, '-42.857142857142')
$ws.Range("H26").Value = -42.857142857142

# I26: ('15', None, '12') -> ('15', None, '14')
$ws.Range("I26").Value = 14

# J26: ('15', None, '17') -> ('15', None, '20')
$ws.Range("J26").Value = 20

# K26: ('16', None, '-29.411764705882') -> ('16', None, '-30')
$ws.Range("K26").Value = -30

# L26: ('16', None, '-14.285714285714') -> ('16', None, '-12.5')
$ws.Range("L26").Value = -12.5

# M26: ('16', None, '0') -> ('16', None, '16.666666666666')
$ws.Range("M26").Value = 16.666666666666

# G28: ('15', None, '1') -> ('15', None, '2')
$ws.Range("G28").Value = 2

# J28: ('15', None, '10') -> ('15', None, '11')
$ws.Range("J28").Value = 11

# K28: ('16', None, '-70') -> ('16', None, '-72.727272727272')
$ws.Range("K28").Value = -72.727272727272

# L28: ('16', None, '0') -> ('16', None, '-25')
$ws.Range("L28").Value = -25

# D31: ('15', None, '1') -> ('14', 's', '20')
$ws.Range("C20").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4122) | Out-Null
$ws.Range("D31").Value = "0"

# E31: ('16', None, '-100') -> ('14', 's', '21')
$ws.Range("C20").Copy() | Out-Null
$ws.Range("E31").PasteSpecial(-4122) | Out-Null
$ws.Range("E31").Value = "***.*"

$excel.CutCopyMode = $false